# Update header labels so that, once imported into Power BI, the first row
# can be automatically promoted to column headers.
# Sheets 1, 2, 3, 5 and 6 use the "Ano" prefix; sheet 4 uses "Intervalo".

$wb = $excel.ActiveWorkbook

$anoSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Emissoes Totais (MtCO2eq)"
)

foreach ($sheetName in $anoSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B1").Value2 = "Ano " + $ws.Range("B1").Value2
    $ws.Range("C1").Value2 = "Ano " + $ws.Range("C1").Value2
    $ws.Range("D1").Value2 = "Ano " + $ws.Range("D1").Value2
    $ws.Range("E1").Value2 = "Ano " + $ws.Range("E1").Value2
}

$wsIntervalo = $wb.Worksheets.Item("Potencia Incremental - SIN(MW)")
$wsIntervalo.Range("B1").Value2 = "Intervalo " + $wsIntervalo.Range("B1").Value2
$wsIntervalo.Range("C1").Value2 = "Intervalo " + $wsIntervalo.Range("C1").Value2
$wsIntervalo.Range("D1").Value2 = "Intervalo " + $wsIntervalo.Range("D1").Value2
$wsIntervalo.Range("E1").Value2 = "Intervalo " + $wsIntervalo.Range("E1").Value2

$wsCusto = $wb.Worksheets.Item("Custo Total (bilhões de R$)")
$wsCusto.Range("B1").Value2 = "Ano " + $wsCusto.Range("B1").Value2
